$d = $word.ActiveDocument

function Add-BitsRun {
    param(
        [string]$FindText,
        [string]$InsertText
    )

    $rng = $d.Content
    $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Collapse(0)

    $insertStart = $rng.Start
    $rng.InsertAfter($InsertText)
    $insertEnd = $rng.End

    # Re-select just the inserted text and nudge a character property on and
    # back off again so Word is forced to materialize it as its own run
    # (matching formatting to the rest of the line) instead of silently
    # merging it into the preceding run.
    $newRng = $d.Range($insertStart, $insertEnd)
    $newRng.Font.Bold = 1
    $newRng.Font.Bold = 0
}

Add-BitsRun "byte: Byte" ":  8 bits"
Add-BitsRun "short: Short" ": 16 bits"
Add-BitsRun "int: Integer" ": 32 bits"
Add-BitsRun "long: Long" ": 64 bits"
